$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trulia")

# Rename the sheet "Trulia" -> "ManageAppointment"
$ws.Name = "ManageAppointment"

# Header row
$ws.Range("A1").Value = "email"
$ws.Range("B1").Value = "code"

# Data row: email address (hyperlinked) + numeric verification code (kept as text)
$ws.Range("A2").Value = "osman@yahoo.com"
[void]$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:osman@yahoo.com")
$ws.Range("A2").Style = "Hyperlink"

$ws.Range("B2").Value = "'736475372"

# Column widths to fit the new content
$ws.Columns.Item(1).ColumnWidth = 21
$ws.Columns.Item(2).ColumnWidth = 11

# Make this the active/selected sheet and cell
[void]$ws.Range("B2").Select()
$ws.Activate()
